$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Apoe"
$ws.Range("C2").Value = "Vldlr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 67.77251700000001
$ws.Range("H2").Value = 203.317551
$ws.Range("I2").Value = 0.4079637943863715
$ws.Range("J2").Value = 0.4079637943863715
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06089466666666667
$ws.Range("N2").Value = 0.182684
$ws.Range("O2").Value = 0.001903591634475228
$ws.Range("P2").Value = 0.001903591634475228
$ws.Range("Q2").Value = 4.126984831876001
$ws.Range("R2").Value = 37.142863486884
$ws.Range("S2").Value = 0.0007765964661626687
$ws.Range("T2").Value = 0.0007765964661626686

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Apoe"
$ws.Range("C3").Value = "Vldlr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 67.77251700000001
$ws.Range("H3").Value = 203.317551
$ws.Range("I3").Value = 0.4079637943863715
$ws.Range("J3").Value = 0.4079637943863715
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 28.046323
$ws.Range("N3").Value = 84.138969
$ws.Range("O3").Value = 0.8767392739472014
$ws.Range("P3").Value = 0.8767392739472013
$ws.Range("Q3").Value = 1900.769902304991
$ws.Range("R3").Value = 17106.92912074492
$ws.Range("S3").Value = 0.3576778808870527
$ws.Range("T3").Value = 0.3576778808870527

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Apoe"
$ws.Range("C4").Value = "Vldlr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 67.77251700000001
$ws.Range("H4").Value = 203.317551
$ws.Range("I4").Value = 0.4079637943863715
$ws.Range("J4").Value = 0.4079637943863715
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.882136333333333
$ws.Range("N4").Value = 11.646409
$ws.Range("O4").Value = 0.1213571344183235
$ws.Range("P4").Value = 0.1213571344183235
$ws.Range("Q4").Value = 263.1021506471511
$ws.Range("R4").Value = 2367.919355824359
$ws.Range("S4").Value = 0.04950931703315616
$ws.Range("T4").Value = 0.04950931703315616

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Apoe"
$ws.Range("C5").Value = "Vldlr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 60.97760633333333
$ws.Range("H5").Value = 182.932819
$ws.Range("I5").Value = 0.3670611149405164
$ws.Range("J5").Value = 0.3670611149405164
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06089466666666667
$ws.Range("N5").Value = 0.182684
$ws.Range("O5").Value = 0.001903591634475228
$ws.Range("P5").Value = 0.001903591634475228
$ws.Range("Q5").Value = 3.713211011799556
$ws.Range("R5").Value = 33.418899106196
$ws.Range("S5").Value = 0.0006987344677419171
$ws.Range("T5").Value = 0.000698734467741917

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Apoe"
$ws.Range("C6").Value = "Vldlr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 60.97760633333333
$ws.Range("H6").Value = 182.932819
$ws.Range("I6").Value = 0.3670611149405164
$ws.Range("J6").Value = 0.3670611149405164
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 28.046323
$ws.Range("N6").Value = 84.138969
$ws.Range("O6").Value = 0.8767392739472014
$ws.Range("P6").Value = 0.8767392739472013
$ws.Range("Q6").Value = 1710.197642991512
$ws.Range("R6").Value = 15391.77878692361
$ws.Range("S6").Value = 0.3218168954071986
$ws.Range("T6").Value = 0.3218168954071985

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Apoe"
$ws.Range("C7").Value = "Vldlr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 60.97760633333333
$ws.Range("H7").Value = 182.932819
$ws.Range("I7").Value = 0.3670611149405164
$ws.Range("J7").Value = 0.3670611149405164
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.882136333333333
$ws.Range("N7").Value = 11.646409
$ws.Range("O7").Value = 0.1213571344183235
$ws.Range("P7").Value = 0.1213571344183235
$ws.Range("Q7").Value = 236.7233810663301
$ws.Range("R7").Value = 2130.510429596971
$ws.Range("S7").Value = 0.04454548506557593
$ws.Range("T7").Value = 0.04454548506557592

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Apoe"
$ws.Range("C8").Value = "Vldlr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 37.37372866666667
$ws.Range("H8").Value = 112.121186
$ws.Range("I8").Value = 0.2249750906731122
$ws.Range("J8").Value = 0.2249750906731122
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.06089466666666667
$ws.Range("N8").Value = 0.182684
$ws.Range("O8").Value = 0.001903591634475228
$ws.Range("P8").Value = 0.001903591634475228
$ws.Range("Q8").Value = 2.275860749247112
$ws.Range("R8").Value = 20.482746743224
$ws.Range("S8").Value = 0.0004282607005706423
$ws.Range("T8").Value = 0.0004282607005706422

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Apoe"
$ws.Range("C9").Value = "Vldlr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 37.37372866666667
$ws.Range("H9").Value = 112.121186
$ws.Range("I9").Value = 0.2249750906731122
$ws.Range("J9").Value = 0.2249750906731122
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 28.046323
$ws.Range("N9").Value = 84.138969
$ws.Range("O9").Value = 0.8767392739472014
$ws.Range("P9").Value = 0.8767392739472013
$ws.Range("Q9").Value = 1048.195665899693
$ws.Range("R9").Value = 9433.760993097236
$ws.Range("S9").Value = 0.1972444976529502
$ws.Range("T9").Value = 0.1972444976529502

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Apoe"
$ws.Range("C10").Value = "Vldlr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 37.37372866666667
$ws.Range("H10").Value = 112.121186
$ws.Range("I10").Value = 0.2249750906731122
$ws.Range("J10").Value = 0.2249750906731122
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.882136333333333
$ws.Range("N10").Value = 11.646409
$ws.Range("O10").Value = 0.1213571344183235
$ws.Range("P10").Value = 0.1213571344183235
$ws.Range("Q10").Value = 145.0899099690082
$ws.Range("R10").Value = 1305.809189721074
$ws.Range("S10").Value = 0.02730233231959139
$ws.Range("T10").Value = 0.02730233231959139
